# Zero-pad the subject tags in column A (S1..S6 -> S01..S06), and update
# the sheet view's topLeftCell / selection to match the post-edit state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2:A9").Value = "S01"
$ws.Range("A10:A17").Value = "S02"
$ws.Range("A18:A25").Value = "S03"
$ws.Range("A26:A33").Value = "S04"
$ws.Range("A34:A41").Value = "S05"
$ws.Range("A42:A49").Value = "S06"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("A44").Select()
